$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-180)
# by bumping the date serial value from 45180 to 45181 (one day later).
$ws.Range("C2:C180").Value = 45181
